$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: reuse the panorama_fish_eye icon instead of
# the now-removed account_tree / all_inbox / account_box icons ---
$ws.Range("E2").Value = "panorama_fish_eye"
$ws.Range("E5").Value = "panorama_fish_eye"
$ws.Range("E8").Value = "panorama_fish_eye"

# --- Add new menu rows 20-26 (Usuarios, Solicitudes, Servicios) ---
$ws.Range("C20").Value = "Usuarios"
$ws.Range("D20").Value = "usuarios"
$ws.Range("E20").Value = "supervisor_account"
$ws.Range("C21").Value = "Usuarios nuevo"
$ws.Range("I21").Value = "Ruta para crear nuevo usuario"
$ws.Range("C22").Value = "Usuarios editar"
$ws.Range("D22").Value = "usuarios/editar"
$ws.Range("D21").Value = "usuarios/nuevo"
$ws.Range("I22").Value = "Ruta para editar un usuario"
$ws.Range("C23").Value = "Solicitudes"
$ws.Range("D23").Value = "solicitudes"
$ws.Range("I23").Value = "Menú para listar todas las solicitudes"
$ws.Range("C24").Value = "Solicitudes nueva"
$ws.Range("D24").Value = "solicitudes/nuevo"
$ws.Range("E23").Value = "assignment"
$ws.Range("I24").Value = "Ruta para nueva solicitud"
$ws.Range("C25").Value = "Solicitudes editar"
$ws.Range("D25").Value = "solicitudes/editar"
$ws.Range("I25").Value = "Ruta para editar solicitud"
$ws.Range("C26").Value = "Servicios"
$ws.Range("D26").Value = "servicios"
$ws.Range("I26").Value = "Menú para listar los servicios"
$ws.Range("I20").Value = "Menú para listar todos los usuarios"
$ws.Range("E26").Value = "style"
$ws.Range("A20").Value = 20
$ws.Range("B20").Value = 0
$ws.Range("F20").Value = "visible"
$ws.Range("G20").Value = "Digitador"
$ws.Range("H20").Value = 3
$ws.Range("A21").Value = 21
$ws.Range("B21").Value = 20
$ws.Range("E21").Value = "minimize"
$ws.Range("F21").Value = "oculto"
$ws.Range("G21").Value = "Digitador"
$ws.Range("H21").Value = 0
$ws.Range("A22").Value = 22
$ws.Range("B22").Value = 20
$ws.Range("E22").Value = "minimize"
$ws.Range("F22").Value = "oculto"
$ws.Range("G22").Value = "Digitador"
$ws.Range("H22").Value = 0
$ws.Range("A23").Value = 23
$ws.Range("B23").Value = 0
$ws.Range("F23").Value = "visible"
$ws.Range("G23").Value = "Digitador"
$ws.Range("H23").Value = 4
$ws.Range("A24").Value = 24
$ws.Range("B24").Value = 23
$ws.Range("E24").Value = "minimize"
$ws.Range("F24").Value = "oculto"
$ws.Range("G24").Value = "Digitador"
$ws.Range("H24").Value = 0
$ws.Range("A25").Value = 25
$ws.Range("B25").Value = 23
$ws.Range("E25").Value = "minimize"
$ws.Range("F25").Value = "oculto"
$ws.Range("G25").Value = "Digitador"
$ws.Range("H25").Value = 0
$ws.Range("A26").Value = 26
$ws.Range("B26").Value = 0
$ws.Range("F26").Value = "visible"
$ws.Range("G26").Value = "Digitador"
$ws.Range("H26").Value = 5

# --- Restore the active cell selection ---
[void]$ws.Range("H21").Select()
